$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Task Description"
$ws.Range("B1").Value = "Start Date"
$ws.Range("C1").Value = "Completion Date"
$ws.Range("D1").Value = "Milestone Completion Date"

# Task list (column A, rows 2-16)
$ws.Range("A2").Value = "Specify Game Mechanics"
$ws.Range("A3").Value = "Develop Game Prototype"
$ws.Range("A4").Value = "Integrate Unity and RL Framework"
$ws.Range("A5").Value = "Design RL Agent Architecture"
$ws.Range("A6").Value = "Implement RL Agent"
$ws.Range("A7").Value = "Create Training Environment"
$ws.Range("A8").Value = "Define Training Scenarios"

# "Evaluate RL Agent" was typed first, then "Train RL Agent" was inserted
# above it, pushing it down a row - matches original authoring order
# reflected in the shared-strings table.
$ws.Range("A9").Value = "Evaluate RL Agent"
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "Train RL Agent"

$ws.Range("A11").Value = "Refine Game Mechanics"
$ws.Range("A12").Value = "Enhance Game Features"
$ws.Range("A13").Value = "Update RL Agent"
$ws.Range("A14").Value = "Test Gameplay Experience"
$ws.Range("A15").Value = "Document Project Progress"
$ws.Range("A16").Value = "Project Presentation and Reporting"

# Column widths to match the target layout
$ws.Columns.Item(1).ColumnWidth = 28.5859375
$ws.Columns.Item(2).ColumnWidth = 8.703125
$ws.Columns.Item(3).ColumnWidth = 14.05859375
$ws.Columns.Item(4).ColumnWidth = 22.46875

$ws.Range("D16").Select()
